$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- values that used to belong to row 4
$ws.Range("D2").Value = 45063
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 21000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21500
$ws.Range("P2").Value = 1433

# Row 4 <- values that used to belong to row 6
$ws.Range("D4").Value = 45084
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 23000
$ws.Range("M4").Value = 22556
$ws.Range("P4").Value = 1504

# Row 6 <- values that used to belong to row 2
$ws.Range("D6").Value = 44749
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17556
$ws.Range("P6").Value = 1170
